$wb = $excel.ActiveWorkbook

# Rename "mySheet" to "GEUK 도서 리스트"
$ws = $wb.Worksheets.Item("mySheet")
$ws.Name = "GEUK 도서 리스트"

# Shift the whole table one column to the right (A:E -> B:F),
# preserving existing formatting / styles / column widths.
$ws.Columns.Item(1).Insert()

# New header row
$ws.Range("B1").Value = "No."
$ws.Range("C1").Value = "Title"
$ws.Range("D1").Value = "Author"
$ws.Range("E1").Value = "Released Date"
$ws.Range("F1").Value = "Note"
$ws.Range("G1").Value = "분실여부"

# New "lost" status column
$ws.Range("G3").Value = $false

# Update the active selection to the new anchor cell
$ws.Range("G4").Select()
